# Add new simulation result rows (110-133) to Sheet1, extending the
# existing data table that currently spans A1:I109 up to A1:I133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10,0.01,0.01,-0.5,485.389730671534,4.622759339728895,1.58613062912976,0,29),
    @(10,0.01,0.01,0,457.3696527705521,4.234904192319927,2.339443072459451,9,38),
    @(10,0.01,0.01,0.5,456.2853832857861,4.517677062235506,1.322501343264093,0,50),
    @(10,0.01,0.5,-0.5,126.4616004450432,1.227782528592652,1.512833952334825,37,40),
    @(10,0.01,0.5,0,301.4240134172564,2.667469145285455,1.811860493372949,11,31),
    @(10,0.01,0.5,0.5,160,1.6,1.113552872566004,17,50),
    @(10,0.01,0.95,-0.5,235.3293326975517,2.353293326975517,2.615723708118364,33,33),
    @(10,0.01,0.95,0,413.2352620545445,3.862011794902285,1.943309980698489,2,28),
    @(10,0.01,0.95,0.5,441.0037963869226,4.410037963869226,1.353353300964745,0,48),
    @(10,0.5,0.01,-0.5,446.4981963964198,4.212247135815281,1.150222501373096,0,30),
    @(10,0.5,0.01,0,275.3303953583768,2.525966929893365,1.60138852172902,21,34),
    @(10,0.5,0.01,0.5,508.2902449861257,4.88740620178967,1.134168215823737,0,47),
    @(10,0.5,0.5,-0.5,225.0418148726914,1.771982794273161,1.540093503791875,27,27),
    @(10,0.5,0.5,0,366.8478360995956,2.565369483213955,1.378911489948525,12,30),
    @(10,0.5,0.5,0.5,570.6887643979227,4.076348317128019,1.772153281911704,3,37),
    @(10,0.5,0.95,-0.5,408.5878704683959,4.005763435964666,1.143153739249074,0,36),
    @(10,0.5,0.95,0,418.5036749515784,3.770303377942148,1.557065213122966,0,40),
    @(10,0.5,0.95,0.5,514.3652568658816,4.898716732056015,1.126605985134598,0,55),
    @(10,0.95,0.01,-0.5,762.8841343678141,7.628841343678141,0.8722268930657734,0,66),
    @(10,0.95,0.01,0,732.5759519670314,7.182117176147367,1.443557942985183,0,82),
    @(10,0.95,0.01,0.5,342.5029183196256,3.425029183196256,0.9215069691836295,0,89),
    @(10,0.95,0.5,-0.5,452.2221081874419,3.932366158151669,0.8059403678215855,0,20),
    @(10,0.95,0.5,0,476.4012685208576,4.253582754650514,0.8535000414577344,0,20),
    @(10,0.95,0.5,0.5,432.3106080000993,3.632862252101675,0.7640840675584952,0,14)
)

$startRow = 110
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
